# Lean UX - Canvas adaptado
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title text in the merged header cell (A2:B2)
$ws.Range("A2").Value = "product Backlog - sprint 1"

# Move the active selection to A3 (matches saved cursor position)
$ws.Range("A3").Select()
